{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (under the \"Impact\"\n// sub-heading) from 6 job-duty style bullets into 4 impact-focused\n// accomplishment bullets, per the commit:\n//   \"Fix Key Achievements to use proper accomplishment statements\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Exact original bullet text (in document order) that this edit replaces.\nconst oldBullets = [\n  \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n];\n\n// New bullet text that should exist in their place, in document order.\nconst newBullets = [\n  \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\u2022 $4.7M savings enabled nonprofit access\",\n  \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n  \"\u2022 Real-time collaboration at national scale\"\n];\n\n// Locate the contiguous run of paragraphs matching oldBullets exactly.\nconst items = paragraphs.items;\nlet startIdx = -1;\nfor (let i = 0; i + oldBullets.length <= items.length; i++) {\n  let match = true;\n  for (let j = 0; j < oldBullets.length; j++) {\n    if (items[i + j].text !== oldBullets[j]) {\n      match = false;\n      break;\n    }\n  }\n  if (match) {\n    startIdx = i;\n    break;\n  }\n}\n\nif (startIdx === -1) {\n  throw new Error(\"Could not locate the Key Achievements bullet block to replace.\");\n}\n\n// Overwrite the first 4 matched paragraphs in place with the new text\n// (Paragraph.text is read-only in this API, so use insertText with the\n// \"Replace\" location to swap the whole-paragraph contents)...\nfor (let j = 0; j < newBullets.length; j++) {\n  items[startIdx + j].insertText(newBullets[j], Word.InsertLocation.replace);\n}\n\n// ...then delete the two now-unneeded trailing paragraphs (indices 4 and 5\n// of the original 6-paragraph block).\nfor (let j = oldBullets.length - 1; j >= newBullets.length; j--) {\n  items[startIdx + j].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Exact original bullet text (in document order) that this edit replaces,\n# under \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\". Single-quoted so that\n# PowerShell does no variable interpolation on the literal '$' in the text.\n$oldBullets = @(\n  '\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations',\n  '\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets',\n  '\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis',\n  '\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%',\n  '\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%',\n  '\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy'\n)\n\n# New bullet text that should exist in their place, in document order.\n$newBullets = @(\n  '\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%',\n  '\u2022 $4.7M savings enabled nonprofit access',\n  '\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations',\n  '\u2022 Real-time collaboration at national scale'\n)\n\nfunction Get-ParaText($idx) {\n    return $d.Paragraphs($idx).Range.Text.TrimEnd(\"`r\")\n}\n\n# Locate the contiguous run of paragraphs matching $oldBullets exactly.\n$count = $d.Paragraphs.Count\n$startIdx = -1\nfor ($i = 1; $i -le ($count - $oldBullets.Count + 1); $i++) {\n    $allMatch = $true\n    for ($j = 0; $j -lt $oldBullets.Count; $j++) {\n        if ((Get-ParaText ($i + $j)) -ne $oldBullets[$j]) {\n            $allMatch = $false\n            break\n        }\n    }\n    if ($allMatch) {\n        $startIdx = $i\n        break\n    }\n}\n\nif ($startIdx -eq -1) {\n    throw \"Could not locate the Key Achievements bullet block to replace.\"\n}\n\n# Overwrite the first 4 matched paragraphs in place with the new text...\nfor ($j = 0; $j -lt $newBullets.Count; $j++) {\n    $d.Paragraphs($startIdx + $j).Range.Text = $newBullets[$j]\n}\n\n# ...then delete the two now-unneeded trailing paragraphs. Delete from the\n# end backwards so earlier indices stay valid, and delete the whole\n# paragraph Range (text + paragraph mark) so no blank paragraph is left\n# behind.\nfor ($j = $oldBullets.Count - 1; $j -ge $newBullets.Count; $j--) {\n    $d.Paragraphs($startIdx + $j).Range.Delete()\n}\n"}
